$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(4186, 4727, 4810, 5039, 5167, 5167, 5294, 5294, 5399, 5399, 5399, 5399, 5399, 5399)

$row = 2
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 3).Value = $val
    $row++
}
